# csvデータ配列早見表.xlsx — "complete renewal accordion kumap24a"
#
# The sheet is a flat lookup table (column A = numeric index, column B =
# matching label) running from row 1 (index 0) to row 70 (index 69).
# This change inserts one new row — "備考" (index 5, between the existing
# "番地"/"備考" header rows and "雑ごみ") — pushing every following row
# down by one, so the table now runs through row 71 (index 70, still
# "古紙衣類_14:30~16:30" as the final label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 6; existing rows 6-70 shift down to 7-71.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "備考"

# Column A is a simple running index (row number - 1); re-stamp it for
# every row that shifted down so it stays sequential (6, 7, 8, ...).
for ($r = 7; $r -le 71; $r++) {
    $ws.Range("A" + $r).Value = $r - 1
}

# Match the saved view/selection state: scrolled down near the bottom of
# the (now longer) table, with the first empty row below it selected.
$ws.Range("A72").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
